$d = $word.ActiveDocument

# Locate the paragraph that starts with "Ви сте учесници" (the intro
# paragraph that currently holds a dozen separately-formatted runs) and
# collapse it down to a single plain run containing the full merged text,
# with the constellation name swapped from "сазвежђа Персеус" to
# " Сазвежђе Бик".

$newText = "Ви сте учесници глобалног посматрачког пројекта, који има за циљ да одреди колико је светлосно загађене у средини у којој живите. Посматрајући звезде унутар  Сазвежђе Бик и упоређујући их са приложеним звезданим картама, посматрачи широм света могу на практичном примеру да увиде колико је светлосно загађење у њиховој средини. Кроз учешће у овом пројекту, допринећете целовитијем сагледавању глобалног проблема."

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Ви сте учесници*") {
        $full = $p.Range
        # Range spanning the paragraph's text but excluding the trailing
        # paragraph mark, so the mark (and its rPr) is left untouched.
        $body = $d.Range($full.Start, $full.End - 1)
        $body.Delete()
        $ins = $d.Range($full.Start, $full.Start)
        $ins.InsertAfter($newText)
        break
    }
}
